# Fruta / hortaliza, semanal
# Insert a new weekly record at row 176 (pushing existing rows 176-200 down to 177-201).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("176:176").Insert()

$ws.Range("A176").Value = 1
$ws.Range("B176").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C176").Value = "Arica y Parinacota"
$ws.Range("D176").Value = 45077
$ws.Range("E176").Value = 15
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100106
$ws.Range("H176").Value = "Oleaginosos"
$ws.Range("I176").Value = 100106002
$ws.Range("J176").Value = "Palta"
$ws.Range("K176").Value = "Hass"
$ws.Range("L176").Value = "Segunda"
$ws.Range("M176").Value = 450
$ws.Range("N176").Value = 24000
$ws.Range("O176").Value = 25000
$ws.Range("P176").Value = 24667
$ws.Range("Q176").Value = "`$/bandeja 10 kilos"
$ws.Range("R176").Value = "Perú"
$ws.Range("S176").Value = 2467
$ws.Range("T176").Value = 10
